# Acta N°3 update:
#  1. Add a new attendee "Brayan García (Tarde)" after "Leyder Vera" in the
#     Asistentes list (same bullet-list formatting as the other names).
#  2. "Logo del Grupo" row: Estado "Pendiente." -> "HECHO (TARDE)".
#  3. "Mínimo 40% de SPMP asignado" row: Estado " Pendiente." (leading blue
#     space run + plain run) -> a single plain run reading "HECHO".
#  4. "Enviar la presentación por correo." row: Estado "Pendiente." -> "HECHO".
#  5. "Citar las diapositivas con IEEE" row: Estado "Please rápido, men." ->
#     "HECHO".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the new attendee paragraph right after "Leyder Vera".
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Leyder Vera")) {
        $r = $p.Range
        $r.Collapse(0)
        # Inserting a paragraph break + text this way makes the new
        # paragraph inherit the source paragraph's numbering / spacing /
        # indentation / justification (same as Word splitting a paragraph).
        $r.InsertAfter("`rBrayan García (Tarde)")
        break
    }
}

# ---------------------------------------------------------------------
# 2) "Logo del Grupo" row -> Estado "Pendiente." becomes "HECHO (TARDE)".
#    This is the first standalone "Pendiente." paragraph in the doc.
#    Replace=1 (wdReplaceOne) so only the first match is touched.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Pendiente.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "HECHO (TARDE)", 1) | Out-Null

# ---------------------------------------------------------------------
# 3) "Mínimo 40% de SPMP asignado" row -> Estado " Pendiente." (a leading
#    single-space run in blue, followed by a "Pendiente." run) becomes a
#    single run reading "HECHO" (keeping the *second* run's formatting,
#    i.e. no blue colour). Drop the leading space character (which is the
#    entire first run) then rewrite the remaining text.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith(" Pendiente.")) {
        $full = $p.Range
        $spaceChar = $d.Range($full.Start, $full.Start + 1)
        $spaceChar.Delete()
        break
    }
}
$d.Content.Find.Execute("Pendiente.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "HECHO", 1) | Out-Null

# ---------------------------------------------------------------------
# 4) "Enviar la presentación por correo." row -> Estado "Pendiente." becomes
#    "HECHO" (leaves the trailing empty run in that paragraph untouched).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Pendiente.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "HECHO", 1) | Out-Null

# ---------------------------------------------------------------------
# 5) "Citar las diapositivas con IEEE" row -> Estado "Please rápido, men."
#    becomes "HECHO".
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Please rápido, men.", $true, $false, $false, $false,
                         $false, $true, 1, $false, "HECHO", 1) | Out-Null
